$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing row formatting down to the two new rows, the same way
# typing a new row under the existing table in Excel inherits formatting.
$ws.Range("A52:F52").Copy()
$ws.Range("A54:F54").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A53:F53").Copy()
$ws.Range("A55:F55").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New day of data: 2025-10-27 (serial 45957), one row per station.
$ws.Range("A54").Value = 45957
$ws.Range("B54").Value = "四方坪站"
$ws.Range("C54").Value = 8252.82
$ws.Range("D54").Value = 6734.9
$ws.Range("E54").Value = 2845.09
$ws.Range("F54").Value = 356

$ws.Range("A55").Value = 45957
$ws.Range("B55").Value = "高岭站"
$ws.Range("C55").Value = 4238.8
$ws.Range("D55").Value = 3623.79
$ws.Range("E55").Value = 1121.63
$ws.Range("F55").Value = 169

# Move the active selection the way Excel would after entering this data.
$ws.Range("J53").Select()

# Configure page setup like the source file (paper size + orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
